$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "unkown" zone entry (row 6: Region=oromiya, Zone=unkown) is being
# dropped per review feedback. Deleting the whole row shifts the rows
# below it up and lets the engine prune the now-unused "unkown" shared
# string from the table on save.
$ws.Rows("6:6").Delete()
